# Update the "Weekly" worksheet to match the target revision:
#  - Several status cells in column J change from "Thu"/"Sun" to "Done" or "Mon"
#  - H16 text is updated (finished message -> modal wording)
#  - H16/H20/H21/H26/H27 get a yellow highlight fill
#  - I26/I27/I33 text updated
#  - J31/J32 gain a "Done" status that didn't exist before
#  - New rows 35-46 are appended with new backlog items

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# --- Column J status updates (Thu/Sun -> Done, or -> Mon) ---
$ws.Range("J9").Value  = "Done"
$ws.Range("J10").Value = "Done"
$ws.Range("J11").Value = "Done"
$ws.Range("J12").Value = "Done"
$ws.Range("J13").Value = "Done"
$ws.Range("J14").Value = "Done"
$ws.Range("J15").Value = "Done"
$ws.Range("J16").Value = "Mon"
$ws.Range("J17").Value = "Done"
$ws.Range("J18").Value = "Done"
$ws.Range("J19").Value = "Done"
$ws.Range("J20").Value = "Mon"
$ws.Range("J22").Value = "Done"
$ws.Range("J23").Value = "Done"
$ws.Range("J24").Value = "Done"
$ws.Range("J25").Value = "Done"

# --- Text updates ---
$ws.Range("H16").Value = "Congratulations, you've finished - modal"
$ws.Range("I26").Value = "LV"
$ws.Range("I27").Value = "MP"
$ws.Range("I33").Value = "MP/CC"

# --- Yellow highlight fill on a handful of task cells ---
$ws.Range("H16").Interior.Color = 65535
$ws.Range("H20").Interior.Color = 65535
$ws.Range("H21").Interior.Color = 65535
$ws.Range("H26").Interior.Color = 65535
$ws.Range("H27").Interior.Color = 65535

# --- New statuses for rows that previously had none ---
$ws.Range("J31").Value = "Done"
$ws.Range("J32").Value = "Done"

# --- New backlog rows appended below the existing table ---
$ws.Range("H35").Value = "Popups"

$ws.Range("H36").Value = "Small game improvements"
$ws.Range("I36").Value = "LV"

$ws.Range("H37").Value = "Music to stop returning on for each level"
$ws.Range("I37").Value = "LV"

$ws.Range("H38").Value = "Increase points/time"
$ws.Range("I38").Value = "LV"

$ws.Range("H39").Value = "Add intro page"
$ws.Range("I39").Value = "MP/CC"

$ws.Range("H40").Value = "Make writeup outline - sections"

$ws.Range("H41").Value = "Make presentation outline - sections"

$ws.Range("H42").Value = "JS error handling for move submission"
$ws.Range("I42").Value = "LV"

$ws.Range("H43").Value = "Improve level names"
$ws.Range("I43").Value = "CC"

$ws.Range("H44").Value = "Lightning"
$ws.Range("I44").Value = "MP"
$ws.Range("J44").Value = "Done"

$ws.Range("H45").Value = "explosion"
$ws.Range("I45").Value = "MP"

$ws.Range("H46").Value = "sinking"
$ws.Range("I46").Value = "MP"

# --- Selection / view state to match the author's last saved cursor ---
$ws.Activate()
$ws.Range("I39").Select()

Write-Host "Weekly sheet updated"
